# Fix matched counts in summary layout.
# For each sheet, update the "Non-matching Rows" (B6), "(Source1) - Matching Rows" (B8),
# and "(Source2) - Matching Rows" (B11) cells with the corrected values.

$wb = $excel.ActiveWorkbook

$updates = @{
    "full"       = @{ B6 = 39;  B8 = 960; B11 = 960 }
    "left"       = @{ B6 = 7;   B8 = 960; B11 = 960 }
    "right"      = @{ B6 = 32;  B8 = 960; B11 = 960 }
    "inner"      = @{ B6 = 0;   B8 = 960; B11 = 960 }
    "complement" = @{ B6 = 39;  B8 = 0;   B11 = 0 }
    "lcomp"      = @{ B6 = 7;   B8 = 0;   B11 = 0 }
    "rcomp"      = @{ B6 = 32;  B8 = 0;   B11 = 0 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $updates[$sheetName]
    $ws.Range("B6").Value = $vals.B6
    $ws.Range("B8").Value = $vals.B8
    $ws.Range("B11").Value = $vals.B11
}
